$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.405.16"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "3.434.09"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("D5").Value = "'575.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'145.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.87%  "
$ws.Range("D7").Value = "3.434.61"
$ws.Range("E7").Value = "  +1.46%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").Value = "4.020.46"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "'28.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.81%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "3.431.82"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "61.514.94"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  +7.16%  "
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").Value = "'395.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.39%  "
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("D25").Value = "'0.995"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "3.574.21"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  +4.25%  "
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "'8.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'1.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.29%  "
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'23.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("D37").Value = "3.463.32"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").Value = "'0.0785"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("D43").Value = "'26.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("E44").Value = "  +3.33%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").Value = "'42.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "2.594.30"
$ws.Range("E49").Value = "  +3.22%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("E51").Value = "  +2.21%  "
